$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("E26").Value = "CONN HEADER 4POS .100 VERT TIN"
$ws.Range("F26").Value = "Molex"
$ws.Range("G26").Value = "22-23-2041"
$ws.Range("H26").Value = "Digikey"
$ws.Range("I26").Value = "WM4202-ND"
$ws.Range("J26").Value = "http://search.digikey.com/scripts/DkSearch/dksus.dll?Detail&name=WM4202-ND"
$ws.Range("K26").Value = 0.415
$ws.Range("L26").Value = 4
$excel.Calculate()
Write-Host "done"
